$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Row 19 - Swine / pigs
$ws.Range("B19").Value = 0.913986017281333
$ws.Range("C19").Value = -0.742098653688219
$ws.Range("D19").Value = -0.913986017281333
$ws.Range("E19").Value = 0.742098653688219
$ws.Range("F19").Value = 0.0995872123477592
$ws.Range("G19").Value = -0.0352770586735377
$ws.Range("H19").Value = -2.85536455243495
$ws.Range("I19").Value = 0.00123036616467417

# Row 20 - Turkeys
$ws.Range("B20").Value = 0.265268303719946
$ws.Range("C20").Value = -0.057958138222745
$ws.Range("D20").Value = -0.265268303719946
$ws.Range("E20").Value = 0.057958138222745
$ws.Range("F20").Value = 0.0109015410894385
$ws.Range("G20").Value = -0.000675535664163159
$ws.Range("H20").Value = -0.203193824366044
$ws.Range("I20").Value = 0.0000362431280735638
